$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Trang chủ (hiển thị sản phẩm nổi bật, tìm kiếm sản phẩm)" as done by Trương Thuận Hòa
# Copy the formatting already used by the other checkmark cells in this column
# so the new checkmark matches the existing look (same cell style as its neighbours).
$ws.Range("D19").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "√"

# Mark "Xử lý giỏ hàng" as done by Phan Hồng Sơn
$ws.Range("E22").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = "√"

# Fill in sub-task detail text for "Hoàn thiện giao diện cho website"
$ws.Range("C27").Value = "Làm ppt"

# Fill in sub-task detail text for "Viết báo cáo và powerpoint"
$ws.Range("C29").Value = "Viết báo cáo"

# Update the active selection to match the author's last click
$ws.Range("C19").Select() | Out-Null
